$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8163891899954194
$ws.Range("C2").Value = 0.1363548090430484
$ws.Range("D2").Value = 0.1124083216962717
$ws.Range("E2").Value = 0.1353499340762969
$ws.Range("F2").Value = 2.131864999578369
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.1908525536553292
$ws.Range("K2").Value = 0.4121628667790844
$ws.Range("L2").Value = 0.2484800205766007
$ws.Range("O2").Value = 5.740659178069961

$ws.Range("B3").Value = 0.7785231929979091
$ws.Range("C3").Value = 0.1356897337931429
$ws.Range("D3").Value = 0.1105904294790534
$ws.Range("E3").Value = 0.1354456548796019
$ws.Range("F3").Value = 2.141136752964243
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.1919772027745168
$ws.Range("K3").Value = 0.3778670620412754
$ws.Range("L3").Value = 0.2445345036909288
$ws.Range("O3").Value = 5.773909290133929

$ws.Range("B4").Value = 0.7555567052493188
$ws.Range("C4").Value = 0.1352854131634018
$ws.Range("D4").Value = 0.1095130407301568
$ws.Range("E4").Value = 0.1355497092653906
$ws.Range("F4").Value = 2.147739748512365
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.1927271816334795
$ws.Range("K4").Value = 0.3568920994117377
$ws.Range("L4").Value = 0.2422055216404431
$ws.Range("O4").Value = 5.796805925111329

$ws.Range("B5").Value = 0.7462696617940878
$ws.Range("C5").Value = 0.135121686754772
$ws.Range("D5").Value = 0.109083814570063
$ws.Range("E5").Value = 0.1356035281034735
$ws.Range("F5").Value = 2.150659589764587
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.1930477719999253
$ws.Range("K5").Value = 0.3483660032229352
$ws.Range("L5").Value = 0.2412800890201794
$ws.Range("O5").Value = 5.806760723277108

$ws.Range("B6").Value = 0.7447319245639505
$ws.Range("C6").Value = 0.1350945634260299
$ws.Range("D6").Value = 0.1090131366978753
$ws.Range("E6").Value = 0.1356131550350099
$ws.Range("F6").Value = 2.15115826857398
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.1931019104365053
$ws.Range("K6").Value = 0.3469515584799012
$ws.Range("L6").Value = 0.2411278532071535
$ws.Range("O6").Value = 5.808451423383787

$ws.Range("B7").Value = 0.7554311644958887
$ws.Range("C7").Value = 0.1352832008621334
$ws.Range("D7").Value = 0.1095072122001213
$ws.Range("E7").Value = 0.1355503888234537
$ws.Range("F7").Value = 2.14777819874611
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.1927314445944752
$ws.Range("K7").Value = 0.3567770262330043
$ws.Range("L7").Value = 0.2421929450362228
$ws.Range("O7").Value = 5.79693765100086

$ws.Range("B8").Value = 0.8032745992380512
$ws.Range("C8").Value = 0.1361246659011002
$ws.Range("D8").Value = 0.1117734942698903
$ws.Range("E8").Value = 0.1353735578866875
$ws.Range("F8").Value = 2.134873135262595
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.1912280087188556
$ws.Range("K8").Value = 0.4003208068206163
$ws.Range("L8").Value = 0.2471002505861719
$ws.Range("O8").Value = 5.75160912590465

$ws.Range("B9").Value = 0.8993173966392192
$ws.Range("C9").Value = 0.1378059466884238
$ws.Range("D9").Value = 0.1165232029738092
$ws.Range("E9").Value = 0.1353848802571491
$ws.Range("F9").Value = 2.116779194887329
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.1887504721540338
$ws.Range("K9").Value = 0.4863472631274988
$ws.Range("L9").Value = 0.2574616504631564
$ws.Range("O9").Value = 5.682390175139659

$ws.Range("B10").Value = 0.9712067489243736
$ws.Range("C10").Value = 0.139059159465198
$ws.Range("D10").Value = 0.1201963663012577
$ws.Range("E10").Value = 0.1356100221079544
$ws.Range("F10").Value = 2.107872960373726
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.1872159293004962
$ws.Range("K10").Value = 0.5499193840907424
$ws.Range("L10").Value = 0.2655193970175418
$ws.Range("O10").Value = 5.643509214706143

$ws.Range("B11").Value = 1.004193762960426
$ws.Range("C11").Value = 0.1396329762811064
$ws.Range("D11").Value = 0.1219066955654426
$ws.Range("E11").Value = 0.1357592229803934
$ws.Range("F11").Value = 2.10477201346977
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.1865796027802169
$ws.Range("K11").Value = 0.5789161876817843
$ws.Range("L11").Value = 0.269280802383264
$ws.Range("O11").Value = 5.628417951334853

$ws.Range("B12").Value = 1.016725329855802
$ws.Range("C12").Value = 0.1398507815750278
$ws.Range("D12").Value = 0.1225599617791318
$ws.Range("E12").Value = 0.1358224212377515
$ws.Range("F12").Value = 2.103734270774765
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.1863475007573214
$ws.Range("K12").Value = 0.5899072071807439
$ws.Range("L12").Value = 0.2707188374277365
$ws.Range("O12").Value = 5.623076264522723

$ws.Range("B13").Value = 1.01402466138137
$ws.Range("C13").Value = 0.1398038508283577
$ws.Range("D13").Value = 0.12241902106922
$ws.Range("E13").Value = 0.1358085128649442
$ws.Range("F13").Value = 2.10395169821301
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.1863970942667592
$ws.Range("K13").Value = 0.5875396342047168
$ws.Range("L13").Value = 0.2704085246294312
$ws.Range("O13").Value = 5.624210104756969

$ws.Range("B14").Value = 1.005223943210154
$ws.Range("C14").Value = 0.1396508850920526
$ws.Range("D14").Value = 0.121960328261693
$ws.Range("E14").Value = 0.135764288239816
$ws.Range("F14").Value = 2.104683902682709
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.1865603301204821
$ws.Range("K14").Value = 0.5798202172128128
$ws.Range("L14").Value = 0.2693988370462108
$ws.Range("O14").Value = 5.627971012484636

$ws.Range("B15").Value = 0.9998384488271199
$ws.Range("C15").Value = 0.1395572553574738
$ws.Range("D15").Value = 0.1216800933217712
$ws.Range("E15").Value = 0.1357380709686531
$ws.Range("F15").Value = 2.105150172813666
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.1866614702342737
$ws.Range("K15").Value = 0.5750932090996912
$ws.Range("L15").Value = 0.2687821513298303
$ws.Range("O15").Value = 5.630323251768203

$ws.Range("B16").Value = 0.9690566230878517
$ws.Range("C16").Value = 0.1390217322228864
$ws.Range("D16").Value = 0.1200853795258894
$ws.Range("E16").Value = 0.1356012104583755
$ws.Range("F16").Value = 2.108094727519884
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.1872587556736214
$ws.Range("K16").Value = 0.548025884423879
$ws.Range("L16").Value = 0.2652754999740523
$ws.Range("O16").Value = 5.644547679711621

$ws.Range("B17").Value = 0.9502452119734812
$ws.Range("C17").Value = 0.1386941440702003
$ws.Range("D17").Value = 0.1191171182201174
$ws.Range("E17").Value = 0.1355292130800052
$ws.Range("F17").Value = 2.11014445369581
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.1876409725558581
$ws.Range("K17").Value = 0.5314403941825958
$ws.Range("L17").Value = 0.2631487607041691
$ws.Range("O17").Value = 5.653938603798679

$ws.Range("B18").Value = 0.939452180469857
$ws.Range("C18").Value = 0.1385060758654078
$ws.Range("D18").Value = 0.1185639111143075
$ws.Range("E18").Value = 0.1354922071839759
$ws.Range("F18").Value = 2.111412881729734
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.1878666263606199
$ws.Range("K18").Value = 0.5219081911350258
$ws.Range("L18").Value = 0.2619345497111283
$ws.Range("O18").Value = 5.659584365652876

$ws.Range("B19").Value = 0.9358024708086248
$ws.Range("C19").Value = 0.1384424603107632
$ws.Range("D19").Value = 0.1183772439984665
$ws.Range("E19").Value = 0.1354804350614316
$ws.Range("F19").Value = 2.111857723036636
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.1879440277535629
$ws.Range("K19").Value = 0.5186820282750659
$ws.Range("L19").Value = 0.2615249937675372
$ws.Range("O19").Value = 5.661537898328874

$ws.Range("B20").Value = 0.9522449512081721
$ws.Range("C20").Value = 0.1387289801362215
$ws.Range("D20").Value = 0.1192198076490882
$ws.Range("E20").Value = 0.1355364216177932
$ws.Range("F20").Value = 2.109916997244909
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.1875996834465177
$ws.Range("K20").Value = 0.5332051933440312
$ws.Range("L20").Value = 0.2633742216803512
$ws.Range("O20").Value = 5.652913636230409

$ws.Range("B21").Value = 1.007807845237124
$ws.Range("C21").Value = 0.1396958010822829
$ws.Range("D21").Value = 0.1220949059279945
$ws.Range("E21").Value = 0.1357770964961134
$ws.Range("F21").Value = 2.104465132748146
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.1865121434491606
$ws.Range("K21").Value = 0.5820873152120214
$ws.Range("L21").Value = 0.269695036631262
$ws.Range("O21").Value = 5.626856220217491

$ws.Range("B22").Value = 1.044354791962178
$ws.Range("C22").Value = 0.1403306563405735
$ws.Range("D22").Value = 0.1240065706503941
$ws.Range("E22").Value = 0.1359734260366139
$ws.Range("F22").Value = 2.101697680817722
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.1858530138102381
$ws.Range("K22").Value = 0.6140958834581056
$ws.Range("L22").Value = 0.2739056939459488
$ws.Range("O22").Value = 5.612000486694683

$ws.Range("B23").Value = 1.024827896290589
$ws.Range("C23").Value = 0.1399915565689085
$ws.Range("D23").Value = 0.122983314897624
$ws.Range("E23").Value = 0.1358650788911113
$ws.Range("F23").Value = 2.103101975023009
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.1862000844381839
$ws.Range("K23").Value = 0.5970068954903525
$ws.Range("L23").Value = 0.2716511385668383
$ws.Range("O23").Value = 5.619730401523356

$ws.Range("B24").Value = 0.9513408008293993
$ws.Range("C24").Value = 0.1387132299188849
$ws.Range("D24").Value = 0.1191733709794818
$ws.Range("E24").Value = 0.1355331489729252
$ws.Range("F24").Value = 2.110019549906681
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.1876183318357718
$ws.Range("K24").Value = 0.5324073182047186
$ws.Range("L24").Value = 0.2632722643515422
$ws.Range("O24").Value = 5.653376255292017

$ws.Range("B25").Value = 0.8731002884531733
$ws.Range("C25").Value = 0.1373478857393025
$ws.Range("D25").Value = 0.1152058520584163
$ws.Range("E25").Value = 0.1353436344024033
$ws.Range("F25").Value = 2.120902909898206
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.1893704493167796
$ws.Range("K25").Value = 0.4630086826598756
$ws.Range("L25").Value = 0.2545800857199652
$ws.Range("O25").Value = 5.699011697198699
